# Rename the custom paragraph style "Docdate" -> "docDate".
#
# Word's object model has no supported way to change a style's
# internal styleId once the style exists (Style.NameLocal only
# rewrites the <w:name> display name, leaving w:styleId untouched).
# The only operation that mints a *new* styleId is Styles.Add, so we
# recreate the style under the new name/id and restore the formatting
# the original "Docdate" style carried (based on Title, quick style,
# 16pt / half-points 32 font size) before removing the old definition.

$d = $word.ActiveDocument

$old = $d.Styles("Docdate")
$baseStyleName = $old.BaseStyle.NameLocal
$isQuickStyle = $old.QuickStyle
$fontSize = $old.Font.Size

$old.Delete()

$new = $d.Styles.Add("docDate", 1)
$new.BaseStyle = $baseStyleName
$new.QuickStyle = $isQuickStyle
$new.Font.Size = $fontSize
